$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the absolute Windows paths to CONTCAR files with relative references
$ws.Range("I3").Value = ".\H2O\CONTCAR"
$ws.Range("I4").Value = ".\H2\CONTCAR"
$ws.Range("I5").Value = ".\O2\CONTCAR"

# Update the current selection, as reflected in the saved file
$ws.Range("I6").Select()
